# Update tissue name values on the "Samples" sheet to lowercase / abbreviated
# forms (fixture-style short tissue descriptors), and move the active
# selection on that sheet to E18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

$ws.Range("D3").Value  = "brain"
$ws.Range("D4").Value  = "diaph"
$ws.Range("D5").Value  = "gast"
$ws.Range("D7").Value  = "heart"
$ws.Range("D8").Value  = "kidney"
$ws.Range("D9").Value  = "liver"
$ws.Range("D10").Value = "lung"
$ws.Range("D11").Value = "pancreas"
$ws.Range("D12").Value = "quad"
$ws.Range("D13").Value = "SmIn"
$ws.Range("D14").Value = "soleus"
$ws.Range("D15").Value = "spleen"
$ws.Range("D16").Value = "serum"

$ws.Activate()
$ws.Range("E18").Select()
